# Add a new ballot row (row 16) for voter "Steve Politi" (NJ.com), matching
# the "update with Politi ballot" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ballots")

# Voter name
$ws.Range("A16").Value = "Steve Politi"

# His "x" votes for: Barry Bonds(C), Roger Clemens(D), Roy Halladay(E),
# Edgar Martinez(I), Mike Mussina(K), Andy Pettitte(M), Manny Ramirez(N),
# Mariano Rivera(O), Curt Schilling(Q), Larry Walker(V)
foreach ($col in @("C","D","E","I","K","M","N","O","Q","V")) {
    $ws.Range("$col" + "16").Value = "x"
}

# Vote count, source and date
$ws.Range("AK16").Value = 10
$ws.Range("AL16").Value = "NJ.com"

# Copy the date cell's formatting (short date, numFmtId 14) from the row
# above, then overwrite with the new ballot's date serial (2018-12-05).
$ws.Range("AM15").Copy($ws.Range("AM16"))
$ws.Range("AM16").Value = 43439

# Match the saved selection/active cell from the edit.
[void]$ws.Range("C16").Select()
